$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 2 (shifts existing rows 2-21 down to 3-22)
$ws.Rows.Item(2).Insert()

# Copy the style of row 3 (the row right below, which used to be row 2) onto the new row 2
# so borders/fills/fonts carry over correctly, then set the specific values/styles per the diff.
$ws.Range("A2:E2").Value = $null

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "Git"

# Apply the same style as the similar "no-subitem" rows (e.g. row 12/14 pattern: plain style for A, s=1 style for B/C/D/E)
$ws.Range("B2:E2").Style = $ws.Range("B12").Style

# Update the active selection to match the target state
$ws.Range("B6").Select()
